$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 100; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
